$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'72.363.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.57%  '
$ws.Range("D3").Formula = "'4.051.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.14%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Formula = "'521.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("E6").Value = '  +2.46%  '
$ws.Range("D7").Formula = "'0.716"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +16.85%  '
$ws.Range("D8").Formula = "'4.041.61"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.14%  '
$ws.Range("D9").Formula = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("D10").Formula = "'0.776"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.14%  '
$ws.Range("E11").Value = '  +5.01%  '
$ws.Range("D12").Formula = "'0.0000333"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").Formula = "'48.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +15.52%  '
$ws.Range("D14").Formula = "'11.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.33%  '
$ws.Range("D15").Formula = "'4.695.89"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.13%  '
$ws.Range("D16").Formula = "'4.083.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.07%  '
$ws.Range("D17").Formula = "'21.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.42%  '
$ws.Range("D18").Formula = "'14.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("D21").Formula = "'72.325.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.59%  '
$ws.Range("D22").Formula = "'446.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.44%  '
$ws.Range("D23").Formula = "'104.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +19.67%  '
$ws.Range("E24").Value = '  +6.47%  '
$ws.Range("D25").Formula = "'15.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.75%  '
$ws.Range("E26").Value = '  +1.15%  '
$ws.Range("D27").Formula = "'11.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.06%  '
$ws.Range("D28").Formula = "'11.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.92%  '
$ws.Range("D29").Formula = "'38.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = '  +2.60%  '
$ws.Range("D31").Formula = "'3.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +15.81%  '
$ws.Range("D32").Formula = "'13.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.48%  '
$ws.Range("E33").Value = '  +4.14%  '
$ws.Range("D34").Formula = "'681.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.51%  '
$ws.Range("E35").Value = '  +14.77%  '
$ws.Range("D36").Formula = "'67.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("D37").Formula = "'42.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.67%  '
$ws.Range("D38").Formula = "'0.0₃0878"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.63%  '
$ws.Range("D39").Formula = "'0.432"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.56%  '
$ws.Range("E40").Value = '  +3.40%  '
$ws.Range("E41").Value = '  +6.02%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").Formula = "'0.0502"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.66%  '
$ws.Range("D44").Formula = "'0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Formula = "'3.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("E46").Value = '  +12.15%  '
$ws.Range("E47").Value = '  -1.83%  '
$ws.Range("D48").Formula = "'3.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.89%  '
$ws.Range("D49").Formula = "'9.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +12.62%  '
$ws.Range("D50").Formula = "'3.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.22%  '
$ws.Range("D51").Formula = "'0.000284"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.99%  '

Write-Output "done"